$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (CasesTab) query: drop the trailing "Cohort" output column so the
# query only returns the Response to Treatment as its final field.
$casesTabQuery = "MATCH (s:study)<-[*]-(c:case)<--(demo:demographic)
WHERE s.clinical_study_designation IN ['COTC007B','NCATS-COP01','GLIOMA01']
MATCH (c)<--(diag:diagnosis)
OPTIONAL MATCH (samp:sample)-->(c)
OPTIONAL MATCH (co:cohort)<-[*]-(c)
WITH DISTINCT c, s, demo, diag, co
RETURN  coalesce(c.case_id, '') AS ``Case ID`` ,
        coalesce(s.clinical_study_designation, '') AS ``Study Code`` ,
        coalesce(s.clinical_study_type, '') AS  ``Study Type``,
        coalesce(demo.breed, '') AS Breed ,
        coalesce(diag.disease_term, '') AS Diagnosis ,
        coalesce(diag.stage_of_disease, '') AS ``Stage of Disease`` ,
        coalesce(demo.patient_age_at_enrollment, '') AS Age ,
        coalesce(demo.sex, '') AS Sex ,
        coalesce(demo.neutered_indicator, '') AS ``Neutered Status``,
        coalesce(demo.weight, '') AS ``Weight (kg)``,
        coalesce(diag.best_response, '') AS ``Response to Treatment``"

$ws.Range("B2").Value = $casesTabQuery

# Move the selection to B2 (was D4) and scroll the view back to the top
# (removes the previous topLeftCell="A4" scroll position).
$ws.Range("B2").Select()

$wb.Save()
